# Updates the 'startup' sheet so row 2 becomes the CasesTab query row
# and row 3 becomes the (new) FilesTab query row, with the stat query in
# column C refreshed for both rows (Trials/Cases/Files counts).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text blocks (single-quoted here-strings preserve backticks/quotes/$ literally)
$tabCases = @'
CasesTab
'@
$statQuery = @'
MATCH (s:specimen)-->(c:case)-->(:arm)-->(ct:clinical_trial)
    WHERE c.gender = "FEMALE"
OPTIONAL MATCH (f:file)-->(:sequencing_assay)-->(:nucleic_acid)-->(s)
RETURN 
     COUNT(DISTINCT ct.clinical_trial_designation) AS Trials,
     COUNT(DISTINCT c.case_id) AS Cases,
      COUNT(DISTINCT f) AS Files
'@
$caseQuery = @'
MATCH (c:case)
    WHERE c.gender='FEMALE'
OPTIONAL MATCH (c)-[:of_arm]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f:file)-[*]->(c)
RETURN DISTINCT
    c.case_id AS `Case ID`,
     ct.clinical_trial_designation AS `Trial Code`,
     a.arm_id AS Arm,
      a.arm_drug AS `Arm Treatment`,
c.disease AS Diagnosis,
  c.gender AS Gender,
    c.race AS Race,
    c.ethnicity AS Ethnicity
 
'@
$tabFiles = @'
FilesTab
'@
$fileQuery = @'

MATCH (c:case)
 WHERE c.gender='FEMALE'
   MATCH (f:file)
      OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
       MATCH (f)-[*]->(c)
      OPTIONAL MATCH (f)-->(parent)
      WITH
        f, parent, c, a, ct, 
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, parent, c, a, ct,
        f.file_size /(1024^i) AS value, 
        10^precision AS factor,
        units[i] as unit
WITH    
        f, parent, c, a, ct, unit,
        round(factor * value)/factor AS size
      RETURN DISTINCT 
       f.file_name AS `File Name`,
       head(labels(parent)) as Association,
       f.file_description AS Description,
       f.file_format AS `File Format`,
         CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
       ct.clinical_trial_designation AS `Trial Code`,
       a.arm_id AS Arm,
       c.case_id AS `Case ID`
        
'@

# Row 2 (CasesTab): refresh the stat query before the long case query so new
# unique strings land in the same order the workbook originally recorded them
$ws.Range("A2").Value = $tabCases
$ws.Range("C2").Value = $statQuery
$ws.Range("B2").Value = $caseQuery

# Row 3 (FilesTab, new content replacing the old CasesTab duplicate row)
$ws.Range("A3").Value = $tabFiles
$ws.Range("C3").Value = $statQuery
$ws.Range("B3").Value = $fileQuery

# Row heights grow to fit the longer wrapped query text
$ws.Rows(2).RowHeight = 210
$ws.Rows(3).RowHeight = 409.5

# Leave the active selection on D3, matching the saved view state
$ws.Range("D3").Select() | Out-Null

Write-Output "Updated CasesTab/FilesTab rows on $($ws.Name)"
